$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6174809
$ws.Range("J40").Value = 38463350
$ws.Range("L40").Value = 38463350
$ws.Range("N40").Value = -38463700
$ws.Range("H64").Value = 3846.9285
$ws.Range("I64").Value = 4099.933
$ws.Range("J64").Value = 3214.4167
$ws.Range("K64").Value = 4099.933
$ws.Range("L64").Value = 3214.4167
$ws.Range("M64").Value = -3851.933
$ws.Range("N64").Value = -3710.4167
$ws.Range("H67").Value = 3846.9285
$ws.Range("I67").Value = 4099.933
$ws.Range("J67").Value = 3214.4167
$ws.Range("K67").Value = 4099.933
$ws.Range("L67").Value = 3214.4167
$ws.Range("M67").Value = -3241.933
$ws.Range("N67").Value = -4930.4167
$ws.Range("H74").Value = 3235.0625
$ws.Range("I74").Value = 3397.7896
$ws.Range("J74").Value = 2997.2307
$ws.Range("K74").Value = 3397.7896
$ws.Range("L74").Value = 2997.2307
$ws.Range("M74").Value = -2461.7896
$ws.Range("N74").Value = -4869.2307
$ws.Range("H76").Value = 4172579
$ws.Range("I76").Value = 6417056
$ws.Range("J76").Value = 4264.5713
$ws.Range("K76").Value = 6417056
$ws.Range("L76").Value = 4264.5713
$ws.Range("M76").Value = -6416741
$ws.Range("N76").Value = -4894.5713
$ws.Range("H77").Value = 3235.0625
$ws.Range("I77").Value = 3397.7896
$ws.Range("J77").Value = 2997.2307
$ws.Range("K77").Value = 16988.948
$ws.Range("L77").Value = 14986.1535
$ws.Range("M77").Value = -12308.948
$ws.Range("N77").Value = -24346.1535
$ws.Range("H79").Value = 4172579
$ws.Range("I79").Value = 6417056
$ws.Range("J79").Value = 4264.5713
$ws.Range("K79").Value = 6417056
$ws.Range("L79").Value = 4264.5713
$ws.Range("M79").Value = -6415964
$ws.Range("N79").Value = -6448.5713
$ws.Range("H88").Value = 7851.647
$ws.Range("J88").Value = 9127
$ws.Range("L88").Value = 9127
$ws.Range("N88").Value = -9939
$ws.Range("H91").Value = 7851.647
$ws.Range("J91").Value = 9127
$ws.Range("L91").Value = 9127
$ws.Range("N91").Value = -11935
$ws.Range("H100").Value = 20835784
$ws.Range("I100").Value = 55557092
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 55557092
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -55556551
$ws.Range("N100").Value = -4082
$ws.Range("H125").Value = 2616.6667
$ws.Range("I125").Value = 850
$ws.Range("J125").Value = 3500
$ws.Range("K125").Value = 7650
$ws.Range("L125").Value = 31500
$ws.Range("M125").Value = -5190
$ws.Range("N125").Value = -36420
$ws.Range("H127").Value = 1310.5927
$ws.Range("I127").Value = 522.5
$ws.Range("J127").Value = 1642.421
$ws.Range("K127").Value = 1567.5
$ws.Range("L127").Value = 4927.263
$ws.Range("M127").Value = 3392.5
$ws.Range("N127").Value = -14847.263
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 76925190
$ws.Range("I63").Value = 90910560
$ws.Range("J63").Value = 5699.5
$ws.Range("K63").Value = 90910560
$ws.Range("L63").Value = 5699.5
$ws.Range("M63").Value = -90909874
$ws.Range("N63").Value = -7071.5
$ws.Range("H66").Value = 76925190
$ws.Range("I66").Value = 90910560
$ws.Range("J66").Value = 5699.5
$ws.Range("K66").Value = 454552800
$ws.Range("L66").Value = 28497.5
$ws.Range("M66").Value = -454549368
$ws.Range("N66").Value = -35361.5
$ws.Range("H102").Value = 7409007
$ws.Range("I102").Value = 7409007
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 7409007
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -7407385
$ws.Range("N102").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11918
$ws.Range("I105").Value = 18302.23
$ws.Range("J105").Value = 3618.5
$ws.Range("K105").Value = 18302.23
$ws.Range("L105").Value = 3618.5
$ws.Range("M105").Value = -16555.23
$ws.Range("N105").Value = -7112.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 419
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 419
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H62").Value = 4569.5884
$ws.Range("I62").Value = 4680.364
$ws.Range("J62").Value = 4366.5
$ws.Range("K62").Value = 4680.364
$ws.Range("L62").Value = 4366.5
$ws.Range("M62").Value = -4056.364
$ws.Range("N62").Value = -5614.5
$ws.Range("H65").Value = 4569.5884
$ws.Range("I65").Value = 4680.364
$ws.Range("J65").Value = 4366.5
$ws.Range("K65").Value = 23401.82
$ws.Range("L65").Value = 21832.5
$ws.Range("M65").Value = -20281.82
$ws.Range("N65").Value = -28072.5
$ws.Range("H86").Value = 62502230
$ws.Range("I86").Value = 100002184
$ws.Range("J86").Value = 2316.6667
$ws.Range("K86").Value = 100002184
$ws.Range("L86").Value = 2316.6667
$ws.Range("M86").Value = -100001061
$ws.Range("N86").Value = -4562.6667
$ws.Range("H89").Value = 62502230
$ws.Range("I89").Value = 100002184
$ws.Range("J89").Value = 2316.6667
$ws.Range("K89").Value = 500010920
$ws.Range("L89").Value = 11583.3335
$ws.Range("M89").Value = -500005304
$ws.Range("N89").Value = -22815.3335
$ws.Range("H134").Value = 3270.9092
$ws.Range("I134").Value = 3293.6296
$ws.Range("J134").Value = 3168.6667
$ws.Range("K134").Value = 9880.888800000001
$ws.Range("L134").Value = 9506.000100000001
$ws.Range("M134").Value = -7345.888800000001
$ws.Range("N134").Value = -14576.0001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1031.6
$ws.Range("I2").Value = 1703.5
$ws.Range("J2").Value = 23.75
$ws.Range("K2").Value = 10221
$ws.Range("L2").Value = 142.5
$ws.Range("M2").Value = -10108
$ws.Range("N2").Value = -368.5
$ws.Range("H38").Value = 4348108
$ws.Range("I38").Value = 5882660
$ws.Range("J38").Value = 209.83333
$ws.Range("K38").Value = 17647980
$ws.Range("L38").Value = 629.49999
$ws.Range("M38").Value = -17647633
$ws.Range("N38").Value = -1323.49999
$ws.Range("H131").Value = 2565276
$ws.Range("J131").Value = 1328.75
$ws.Range("L131").Value = 3986.25
$ws.Range("N131").Value = -14066.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5502.5083
$ws.Range("I70").Value = 5554.804
$ws.Range("J70").Value = 5342.1333
$ws.Range("K70").Value = 5554.804
$ws.Range("L70").Value = 5342.1333
$ws.Range("M70").Value = -5284.804
$ws.Range("N70").Value = -5882.1333
$ws.Range("H73").Value = 5502.5083
$ws.Range("I73").Value = 5554.804
$ws.Range("J73").Value = 5342.1333
$ws.Range("K73").Value = 5554.804
$ws.Range("L73").Value = 5342.1333
$ws.Range("M73").Value = -4618.804
$ws.Range("N73").Value = -7214.1333
$ws.Range("H123").Value = 17329.695
$ws.Range("J123").Value = 17567.686
$ws.Range("L123").Value = 17567.686
$ws.Range("N123").Value = -22467.686
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4484.615
$ws.Range("I62").Value = 4227.273
$ws.Range("J62").Value = 5900
$ws.Range("K62").Value = 4227.273
$ws.Range("L62").Value = 5900
$ws.Range("M62").Value = -3603.273
$ws.Range("N62").Value = -7148
$ws.Range("H65").Value = 4484.615
$ws.Range("I65").Value = 4227.273
$ws.Range("J65").Value = 5900
$ws.Range("K65").Value = 21136.365
$ws.Range("L65").Value = 29500
$ws.Range("M65").Value = -18016.365
$ws.Range("N65").Value = -35740
$ws.Range("H126").Value = 956.26086
$ws.Range("I126").Value = 705.26666
$ws.Range("J126").Value = 1426.875
$ws.Range("K126").Value = 2115.79998
$ws.Range("L126").Value = 4280.625
$ws.Range("M126").Value = 354.2000200000002
$ws.Range("N126").Value = -9220.625
$ws.Range("H136").Value = 2814.7454
$ws.Range("I136").Value = 2946.342
$ws.Range("J136").Value = 2520.5881
$ws.Range("K136").Value = 8839.026
$ws.Range("L136").Value = 7561.7643
$ws.Range("M136").Value = -6289.026
$ws.Range("N136").Value = -12661.7643

Write-Host "Applied 210 cell updates across 8 sheets"